$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header-row cells (row 1): field comments for the create/update
#     user + time columns (mirrors the existing comment.* header cells). ---
$ws.Cells.Item(1, 8).Value  = '<%=comment.create_usr_id_lbl%><%selectList.create_usr_id = data.findAllUsr.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.create_usr_id.join(",") }"` })%>'
$ws.Cells.Item(1, 9).Value  = '<%=comment.create_time_lbl%>'
$ws.Cells.Item(1, 10).Value = '<%=comment.update_usr_id_lbl%><%selectList.update_usr_id = data.findAllUsr.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.update_usr_id.join(",") }"` })%>'
$ws.Cells.Item(1, 11).Value = '<%=comment.update_time_lbl%>'

# --- New data-row cells (row 2): model values for the create/update
#     user + time columns (mirrors the existing model.* data cells). ---
$ws.Cells.Item(2, 8).Value  = '<%=model.create_usr_id_lbl%>'
$ws.Cells.Item(2, 9).Value  = '<%~model.create_time ? new Date(model.create_time) : ""%>'
$ws.Cells.Item(2, 10).Value = '<%=model.update_usr_id_lbl%>'
$ws.Cells.Item(2, 11).Value = '<%~model.update_time ? new Date(model.update_time) : ""%>'
